# Generate Report for Handoff
#
# The "b" source file (b.md) has now been handed off for localization, so
# its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", and a new handoff file + handoff datetime is
# recorded for both the zh-cn and de-de locales. The Overview sheet rolls
# this up into the "Latest Handoff Date" column for the b.md row.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet: row 3 = b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = "2016-26-17 14:26:20"

# --- zh-cn sheet: row 3 = b file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-17 14:26:16"

# Keep the D3 hyperlink's displayed text in sync with the new handoff file
# name (the hyperlink target itself is unchanged).
$zhLink = $wsZhCn.Hyperlinks.Item(8)
$zhLink.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# --- de-de sheet: row 3 = b file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-17 14:26:20"

# Keep the D3 hyperlink's displayed text in sync with the new handoff file
# name (the hyperlink target itself is unchanged).
$deLink = $wsDeDe.Hyperlinks.Item(8)
$deLink.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
